# Backlog.xlsx edit: add two new backlog entries (month filter + negative
# balance guard) to the Inventory System Backlog sheet, and move the
# selection onto the newly-added "Status" cell, matching the author's
# commit: "Added month filtering, but need to work on ordered by entry
# date instead of ID, and need to add all entries".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 was a blank feature row: give it a Feature description, make it
# tall enough to show the full text, and mark its Status as "Pending" by
# copying the formatting already used for "Pending" rows (e.g. I13).
$ws.Rows.Item(14).RowHeight = 77
$ws.Range("A14").Value = "Make sure we don't allow negative entry and to have the balance to negative"

$ws.Range("I13:K13").Copy()
$ws.Range("I14:K14").PasteSpecial(-4122)
$ws.Range("I14").Value = "Pending"

# Row 15 was also blank: second new backlog item, also "Pending".
$ws.Range("A15").Value = "Get entries by month filter"

$ws.Range("I13:K13").Copy()
$ws.Range("I15:K15").PasteSpecial(-4122)
$ws.Range("I15").Value = "Pending"

$excel.CutCopyMode = 0

# Reflect where the author's selection ended up after adding the entries.
$null = $ws.Range("I15:K15").Select()
